# Append 5 new "daily export" rows (2025-11-26 .. 2025-11-30) to the
# bottom of the "Chart" sheet's breadcrumb table, mirroring the existing
# rows (Date text in col A, Invalid count in col B, Valid count in col C).
#
# The date strings must land as literal text (same as every other date
# already in column A) rather than being auto-parsed into Excel date
# serials, so each cell's number format is forced to Text ("@") right
# before the assignment and then cleared again afterwards so the cell
# ends up back on the workbook's default (unstyled) format - exactly
# like its neighbours above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 1 }

$newDates = @("2025-11-26", "2025-11-27", "2025-11-28", "2025-11-29", "2025-11-30")
$invalidCount = 0
$validCount = 27

for ($i = 0; $i -lt $newDates.Count; $i++) {
    $row = $lastRow + 1 + $i

    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $newDates[$i]
    $dateCell.ClearFormats()

    $ws.Cells.Item($row, 2).Value = $invalidCount
    $ws.Cells.Item($row, 3).Value = $validCount
}
